$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.521.84"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.619.45"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.12"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.71"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.264"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0613"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.847.59"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "1.620.64"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.99"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "27.521.57"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.85"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.15"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  +6.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.77"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.81"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "1.442.83"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.935"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.562"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.95"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("D47").Value = "1.759.14"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.22"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0997"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.70%  "
